$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.938.33"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "1.663.48"
$ws.Range("E3").Value = "  +1.90%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "217.48"
$ws.Range("E5").Value = "  +1.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.520"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "28.78"
$ws.Range("E8").Value = "  -2.01%  "

$ws.Range("D9").Value = "0.264"
$ws.Range("E9").Value = "  +1.64%  "

$ws.Range("D10").Value = "0.0614"
$ws.Range("E10").Value = "  +0.09%  "

$ws.Range("D11").Value = "0.0901"
$ws.Range("E11").Value = "  -1.76%  "

$ws.Range("D12").Value = "1.901.13"
$ws.Range("E12").Value = "  +1.87%  "

$ws.Range("D13").Value = "1.719.13"
$ws.Range("E13").Value = "  +5.33%  "

$ws.Range("D14").Value = "0.611"
$ws.Range("E14").Value = "  +7.07%  "

$ws.Range("D15").Value = "10.14"
$ws.Range("E15").Value = "  +12.45%  "

$ws.Range("E16").Value = "  +0.65%  "

$ws.Range("D17").Value = "29.948.51"
$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").Value = "65.08"
$ws.Range("E18").Value = "  +0.09%  "

$ws.Range("D19").Value = "241.52"
$ws.Range("E19").Value = "  -1.59%  "

$ws.Range("D20").Value = "0.0₃0712"
$ws.Range("E20").Value = "  +0.74%  "

$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "4.24"
$ws.Range("E22").Value = "  +2.28%  "

$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.00"
$ws.Range("E23").Value = "  +3.83%  "

$ws.Range("D24").Value = "2.17"
$ws.Range("E24").Value = "  +2.25%  "

$ws.Range("D25").Value = "158.05"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").Value = "15.76"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.110"
$ws.Range("E27").Value = "  -1.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.70"
$ws.Range("E28").Value = "  +0.87%  "

$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").Value = "0.0498"
$ws.Range("E30").Value = "  +1.51%  "

$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("E32").Value = "  +1.23%  "

$ws.Range("D33").Value = "3.21"
$ws.Range("E33").Value = "  -0.30%  "

$ws.Range("D34").Value = "1.439.71"
$ws.Range("E34").Value = "  +0.70%  "

$ws.Range("D35").Value = "1.72"
$ws.Range("E35").Value = "  +4.28%  "

$ws.Range("D36").Value = "1.02"
$ws.Range("E36").Value = "  -1.53%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0177"
$ws.Range("E37").Value = "  +3.46%  "

$ws.Range("B38").Value = "Aave"
$ws.Range("C38").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D38").Value = "79.15"
$ws.Range("E38").Value = "  +10.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.580"
$ws.Range("E39").Value = "  +4.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.30"
$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("E41").Value = "  -8.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.850"
$ws.Range("E42").Value = "  +1.75%  "

$ws.Range("D43").Value = "0.0502"
$ws.Range("E43").Value = "  +0.22%  "

$ws.Range("D44").Value = "1.95"
$ws.Range("E44").Value = "  -0.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("D46").Value = "1.01"
$ws.Range("E46").Value = "  -0.38%  "

$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.808.10"
$ws.Range("E47").Value = "  +1.90%  "

$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "50.47"
$ws.Range("E48").Value = "  -7.80%  "

$ws.Range("E49").Value = "  -1.22%  "

$ws.Range("D50").Value = "94.87"
$ws.Range("E50").Value = "  +6.22%  "

$ws.Range("D51").Value = "0.0₆0117"
$ws.Range("E51").Value = "  +7.91%  "
